$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 572, shifting rows 572:613 down to 573:614
$ws.Rows.Item(572).EntireRow.Insert()

# Populate the newly inserted row 572 with the new record.
# Force column A to Text format first so the date-like string "2026/01/05"
# is stored as literal text rather than being auto-converted to a date serial.
$ws.Cells.Item(572, 1).NumberFormat = "@"
$ws.Cells.Item(572, 1).Value = "2026/01/05"
$ws.Cells.Item(572, 2).Value = "月"
$ws.Cells.Item(572, 3).Value = 19
$ws.Cells.Item(572, 4).Value = 201
